# Updates the crypto price/volume table on Sheet1 (rows 2-51) to the latest
# scraped values. For "Price" (column D) entries that look like a plain
# number (e.g. "211.88"), we assign them through a text-literal formula and
# then convert that formula to a static value via Copy/PasteSpecial(values),
# which keeps the cell stored as text (matching the source data) without
# Excel silently reinterpreting it as a number and without leaving behind
# any unused/extra cell styles. Plain non-numeric-looking text (coin names,
# links, and the "Volume(1h)" percentages, which include spaces/%) can be
# assigned directly since Excel keeps those as text automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.495.33'
$ws.Range("E2").Value = '  -0.03%  '

$ws.Range("D3").Value = '1.567.54'
$ws.Range("E3").Value = '  -2.17%  '

$ws.Range("D5").Formula = '="211.88"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  -1.53%  '

$ws.Range("E6").Value = '  -1.35%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").Formula = '="46.16"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  +4.58%  '

$ws.Range("E9").Value = '  +0.01%  '

$ws.Range("E10").Value = '  -1.93%  '

$ws.Range("E11").Value = '  -1.62%  '

$ws.Range("D12").Formula = '="0.0886"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  -0.46%  '

$ws.Range("D13").Value = '1.790.51'
$ws.Range("E13").Value = '  -2.23%  '

$ws.Range("D14").Value = '1.565.54'
$ws.Range("E14").Value = '  -2.33%  '

$ws.Range("E15").Value = '  -2.80%  '

$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").Formula = '="3.68"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = '  -3.20%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '28.493.10'
$ws.Range("E17").Value = '  -0.04%  '

$ws.Range("D18").Formula = '="62.24"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  -1.83%  '

$ws.Range("D19").Formula = '="227.95"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  -2.07%  '

$ws.Range("D20").Formula = '="7.37"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  -2.58%  '

$ws.Range("E21").Value = '  -2.67%  '

$ws.Range("E22").Value = '  +0.07%  '

$ws.Range("E23").Value = '  -6.00%  '

$ws.Range("D24").Formula = '="9.13"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  -3.22%  '

$ws.Range("E25").Value = '  +6.21%  '

$ws.Range("D26").Formula = '="150.87"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  -1.08%  '

$ws.Range("D27").Formula = '="15.01"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  -2.18%  '

$ws.Range("D28").Formula = '="6.46"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  -2.76%  '

$ws.Range("E29").Value = '  -3.98%  '

$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("D31").Formula = '="0.0468"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  -1.84%  '

$ws.Range("D32").Formula = '="1.11"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  -4.07%  '

$ws.Range("E33").Value = '  -1.04%  '

$ws.Range("E34").Value = '  -2.86%  '

$ws.Range("D35").Value = '1.394.25'
$ws.Range("E35").Value = '  -2.18%  '

$ws.Range("E36").Value = '  -0.76%  '

$ws.Range("E37").Value = '  -3.53%  '

$ws.Range("E38").Value = '  +0.98%  '

$ws.Range("D39").Formula = '="2.58"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  +2.40%  '

$ws.Range("E40").Value = '  -1.00%  '

$ws.Range("D41").Formula = '="0.536"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  -1.59%  '

$ws.Range("E42").Value = '  +0.04%  '

$ws.Range("D43").Formula = '="1.89"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  +2.15%  '

$ws.Range("E44").Value = '  -4.56%  '

$ws.Range("E45").Value = '  -4.21%  '

$ws.Range("D46").Formula = '="0.973"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  -0.78%  '

$ws.Range("D47").Formula = '="62.89"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  -3.24%  '

$ws.Range("D48").Value = '1.703.45'
$ws.Range("E48").Value = '  -2.25%  '

$ws.Range("D49").Formula = '="85.95"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  -1.90%  '

$ws.Range("E50").Value = '  -0.29%  '

$ws.Range("E51").Value = '  -4.54%  '
